$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update table data (rows 2-4) with the new schema/database/table names
$ws.Range("A2").Value = "my_database"
$ws.Range("B2").Value = "schema1"
$ws.Range("C2").Value = "kpi_report_q1"

$ws.Range("A3").Value = "data_warehouse"
$ws.Range("B3").Value = "schema2"
$ws.Range("C3").Value = "users"

$ws.Range("A4").Value = "my_database"
$ws.Range("B4").Value = "schema2"
$ws.Range("C4").Value = "kpi_attributes"

# Move the active selection to C9, matching the saved view state
$ws.Range("C9").Select()
